$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.600.75'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.638.96'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.37'
$ws.Range("E5").Value = '  +2.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.90'
$ws.Range("E6").Value = '  +2.02%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.54'
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("E11").Value = '  +2.05%  '
$ws.Range("E12").Value = '  +3.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.094.18'
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.25'
$ws.Range("E14").Value = '  +12.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.588.68'
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.631.27'
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.55'
$ws.Range("E18").Value = '  +2.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.71'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '348.54'
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("E21").Value = '  -1.08%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  -0.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.95'
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.14'
$ws.Range("E27").Value = '  +6.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.04'
$ws.Range("E28").Value = '  +11.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0799'
$ws.Range("E30").Value = '  +5.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.68'
$ws.Range("E31").Value = '  +4.75%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.62'
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.06'
$ws.Range("E34").Value = '  +10.00%  '
$ws.Range("E35").Value = '  +5.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.31'
$ws.Range("E36").Value = '  +8.05%  '
$ws.Range("E37").Value = '  +2.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '333.48'
$ws.Range("E38").Value = '  +13.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.00'
$ws.Range("E39").Value = '  +4.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.55'
$ws.Range("E40").Value = '  +1.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.876'
$ws.Range("E41").Value = '  +3.14%  '
$ws.Range("E42").Value = '  +7.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.70'
$ws.Range("E43").Value = '  +3.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0998'
$ws.Range("E44").Value = '  +1.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '133.02'
$ws.Range("E45").Value = '  -3.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '20.12'
$ws.Range("E46").Value = '  +2.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0558'
$ws.Range("E48").Value = '  +1.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.610'
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("E50").Value = '  +2.27%  '
$ws.Range("E51").Value = '  +0.47%  '
